# Auto-generated edit script: update cached market-price columns (H-N)
# across all 8 leve-profit tables per the scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 11421.286  # H32: was 11995.8
$ws.Cells.Item(32, 10).Value = 11989.8  # J32: was 13326.333
$ws.Cells.Item(32, 12).Value = 11989.8  # L32: was 13326.333
$ws.Cells.Item(32, 14).Value = -12641.8  # N32: was -13978.333
$ws.Cells.Item(64, 8).Value = 11377.667  # H64: was 10739.9
$ws.Cells.Item(64, 10).Value = 13428.571  # J64: was 12375
$ws.Cells.Item(64, 12).Value = 13428.571  # L64: was 12375
$ws.Cells.Item(64, 14).Value = -13924.571  # N64: was -12871
$ws.Cells.Item(67, 8).Value = 11377.667  # H67: was 10739.9
$ws.Cells.Item(67, 10).Value = 13428.571  # J67: was 12375
$ws.Cells.Item(67, 12).Value = 13428.571  # L67: was 12375
$ws.Cells.Item(67, 14).Value = -15144.571  # N67: was -14091
$ws.Cells.Item(74, 8).Value = 6717.6665  # H74: was 7758.857
$ws.Cells.Item(74, 9).Value = 3956.889  # I74: was 3958.2222
$ws.Cells.Item(74, 10).Value = 15000  # J74: was 14600
$ws.Cells.Item(74, 11).Value = 3956.889  # K74: was 3958.2222
$ws.Cells.Item(74, 12).Value = 15000  # L74: was 14600
$ws.Cells.Item(74, 13).Value = -3020.889  # M74: was -3022.2222
$ws.Cells.Item(74, 14).Value = -16872  # N74: was -16472
$ws.Cells.Item(77, 8).Value = 6717.6665  # H77: was 7758.857
$ws.Cells.Item(77, 9).Value = 3956.889  # I77: was 3958.2222
$ws.Cells.Item(77, 10).Value = 15000  # J77: was 14600
$ws.Cells.Item(77, 11).Value = 19784.445  # K77: was 19791.111
$ws.Cells.Item(77, 12).Value = 75000  # L77: was 73000
$ws.Cells.Item(77, 13).Value = -15104.445  # M77: was -15111.111
$ws.Cells.Item(77, 14).Value = -84360  # N77: was -82360
$ws.Cells.Item(98, 8).Value = 2777.3076  # H98: was 2459.0667
$ws.Cells.Item(98, 9).Value = 2777.3076  # I98: was 2459.0667
$ws.Cells.Item(98, 11).Value = 2777.3076  # K98: was 2459.0667
$ws.Cells.Item(98, 13).Value = -1279.3076  # M98: was -961.0666999999999
$ws.Cells.Item(116, 8).Value = 3999.625  # H116: was 4066.5
$ws.Cells.Item(116, 9).Value = 3679.4  # I116: was 3599.6667
$ws.Cells.Item(116, 11).Value = 3679.4  # K116: was 3599.6667
$ws.Cells.Item(116, 13).Value = -237.4000000000001  # M116: was -157.6667000000002
$ws.Cells.Item(122, 8).Value = 2777.3076  # H122: was 2459.0667
$ws.Cells.Item(122, 9).Value = 2777.3076  # I122: was 2459.0667
$ws.Cells.Item(122, 11).Value = 8331.9228  # K122: was 7377.2001
$ws.Cells.Item(122, 13).Value = -5881.9228  # M122: was -4927.2001
$ws.Cells.Item(136, 8).Value = 299998.8  # H136: was 283330.66
$ws.Cells.Item(136, 10).Value = 299998.8  # J136: was 283330.66
$ws.Cells.Item(136, 12).Value = 299998.8  # L136: was 283330.66
$ws.Cells.Item(136, 14).Value = -310198.8  # N136: was -293530.66
$ws.Cells.Item(137, 8).Value = 1631.8966  # H137: was 1660.6
$ws.Cells.Item(137, 9).Value = 1391.8422  # I137: was 1446.9
$ws.Cells.Item(137, 11).Value = 4175.5266  # K137: was 4340.700000000001
$ws.Cells.Item(137, 13).Value = -1625.5266  # M137: was -1790.700000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 1863.375  # H26: was 1488.25
$ws.Cells.Item(26, 9).Value = 1917.8334  # I26: was 1501.4
$ws.Cells.Item(26, 10).Value = 1700  # J26: was 1466.3334
$ws.Cells.Item(26, 11).Value = 1917.8334  # K26: was 1501.4
$ws.Cells.Item(26, 12).Value = 1700  # L26: was 1466.3334
$ws.Cells.Item(26, 13).Value = -1587.8334  # M26: was -1171.4
$ws.Cells.Item(26, 14).Value = -2360  # N26: was -2126.3334
$ws.Cells.Item(32, 8).Value = 2975.913  # H32: was 2971.4639
$ws.Cells.Item(32, 9).Value = 3125.0923  # I32: was 3120.3691
$ws.Cells.Item(32, 11).Value = 3125.0923  # K32: was 3120.3691
$ws.Cells.Item(32, 13).Value = -2838.0923  # M32: was -2833.3691

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1923.8334  # H86: was 1961.8334
$ws.Cells.Item(86, 9).Value = 2010.75  # I86: was 1917.75
$ws.Cells.Item(86, 10).Value = 1750  # J86: was 2050
$ws.Cells.Item(86, 11).Value = 2010.75  # K86: was 1917.75
$ws.Cells.Item(86, 12).Value = 1750  # L86: was 2050
$ws.Cells.Item(86, 13).Value = -887.75  # M86: was -794.75
$ws.Cells.Item(86, 14).Value = -3996  # N86: was -4296
$ws.Cells.Item(89, 8).Value = 1923.8334  # H89: was 1961.8334
$ws.Cells.Item(89, 9).Value = 2010.75  # I89: was 1917.75
$ws.Cells.Item(89, 10).Value = 1750  # J89: was 2050
$ws.Cells.Item(89, 11).Value = 10053.75  # K89: was 9588.75
$ws.Cells.Item(89, 12).Value = 8750  # L89: was 10250
$ws.Cells.Item(89, 13).Value = -4437.75  # M89: was -3972.75
$ws.Cells.Item(89, 14).Value = -19982  # N89: was -21482
$ws.Cells.Item(94, 8).Value = 2792.5833  # H94: was 2826.9167
$ws.Cells.Item(94, 9).Value = 2631.1  # I94: was 2672.3
$ws.Cells.Item(94, 11).Value = 2631.1  # K94: was 2672.3
$ws.Cells.Item(94, 13).Value = -2180.1  # M94: was -2221.3
$ws.Cells.Item(99, 8).Value = 3795.1304  # H99: was 3459.56
$ws.Cells.Item(99, 9).Value = 2607.6155  # I99: was 2406.6667
$ws.Cells.Item(99, 10).Value = 5338.9  # J99: was 5038.9
$ws.Cells.Item(99, 11).Value = 2607.6155  # K99: was 2406.6667
$ws.Cells.Item(99, 12).Value = 5338.9  # L99: was 5038.9
$ws.Cells.Item(99, 13).Value = -1109.6155  # M99: was -908.6667000000002
$ws.Cells.Item(99, 14).Value = -8334.9  # N99: was -8034.9
$ws.Cells.Item(107, 8).Value = 4249.875  # H107: was 4756.857
$ws.Cells.Item(107, 9).Value = 3416.5  # I107: was 3549.6667
$ws.Cells.Item(107, 10).Value = 6750  # J107: was 12000
$ws.Cells.Item(107, 11).Value = 3416.5  # K107: was 3549.6667
$ws.Cells.Item(107, 12).Value = 6750  # L107: was 12000
$ws.Cells.Item(107, 13).Value = -1496.5  # M107: was -1629.6667
$ws.Cells.Item(107, 14).Value = -10590  # N107: was -15840

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(28, 8).Value = 49723.8  # H28: was 49667.25
$ws.Cells.Item(28, 10).Value = 49723.8  # J28: was 49667.25
$ws.Cells.Item(28, 12).Value = 49723.8  # L28: was 49667.25
$ws.Cells.Item(28, 14).Value = -50213.8  # N28: was -50157.25
$ws.Cells.Item(58, 8).Value = 4305.5674  # H58: was 4630.4707
$ws.Cells.Item(58, 9).Value = 1664.56  # I58: was 1806.5454
$ws.Cells.Item(58, 11).Value = 1664.56  # K58: was 1806.5454
$ws.Cells.Item(58, 13).Value = -1461.56  # M58: was -1603.5454
$ws.Cells.Item(95, 8).Value = 13108.333  # H95: was 18331.25
$ws.Cells.Item(95, 10).Value = 13108.333  # J95: was 18331.25
$ws.Cells.Item(95, 12).Value = 13108.333  # L95: was 18331.25
$ws.Cells.Item(95, 14).Value = -18600.333  # N95: was -23823.25
$ws.Cells.Item(105, 8).Value = 1548.4615  # H105: was 1891.25
$ws.Cells.Item(105, 9).Value = 1503  # I105: was 2171.6667
$ws.Cells.Item(105, 10).Value = 1700  # J105: was 1050
$ws.Cells.Item(105, 11).Value = 1503  # K105: was 2171.6667
$ws.Cells.Item(105, 12).Value = 1700  # L105: was 1050
$ws.Cells.Item(105, 13).Value = 244  # M105: was -424.6667000000002
$ws.Cells.Item(105, 14).Value = -5194  # N105: was -4544
$ws.Cells.Item(134, 8).Value = 5977.684  # H134: was 6369.5293
$ws.Cells.Item(134, 9).Value = 3541.4285  # I134: was 3714.3076
$ws.Cells.Item(134, 10).Value = 12799.2  # J134: was 14999
$ws.Cells.Item(134, 11).Value = 10624.2855  # K134: was 11142.9228
$ws.Cells.Item(134, 12).Value = 38397.60000000001  # L134: was 44997
$ws.Cells.Item(134, 13).Value = -8089.2855  # M134: was -8607.9228
$ws.Cells.Item(134, 14).Value = -43467.60000000001  # N134: was -50067
$ws.Cells.Item(136, 8).Value = 4305.5674  # H136: was 4630.4707
$ws.Cells.Item(136, 9).Value = 1664.56  # I136: was 1806.5454
$ws.Cells.Item(136, 11).Value = 4993.68  # K136: was 5419.6362
$ws.Cells.Item(136, 13).Value = -2443.68  # M136: was -2869.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 273.91666  # H34: was 347.77777
$ws.Cells.Item(34, 9).Value = 129.75  # I34: was 153.33333
$ws.Cells.Item(34, 10).Value = 346  # J34: was 445
$ws.Cells.Item(34, 11).Value = 389.25  # K34: was 459.99999
$ws.Cells.Item(34, 12).Value = 1038  # L34: was 1335
$ws.Cells.Item(34, 13).Value = -305.25  # M34: was -375.99999
$ws.Cells.Item(34, 14).Value = -1206  # N34: was -1503
$ws.Cells.Item(39, 8).Value = 971.4286  # H39: was 1066.6666
$ws.Cells.Item(39, 10).Value = 960  # J39: was 1200
$ws.Cells.Item(39, 12).Value = 2880  # L39: was 3600
$ws.Cells.Item(39, 14).Value = -3468  # N39: was -4188
$ws.Cells.Item(55, 8).Value = 663  # H55: was 722.2727
$ws.Cells.Item(55, 10).Value = 733.6923  # J55: was 837.55554
$ws.Cells.Item(55, 12).Value = 2201.0769  # L55: was 2512.66662
$ws.Cells.Item(55, 14).Value = -2555.0769  # N55: was -2866.66662
$ws.Cells.Item(92, 8).Value = 1000  # H92: was 274.5
$ws.Cells.Item(92, 9).Value = 1000  # I92: was 274.5
$ws.Cells.Item(92, 11).Value = 3000  # K92: was 823.5
$ws.Cells.Item(92, 13).Value = -1752  # M92: was 424.5
$ws.Cells.Item(106, 8).Value = 3406.25  # H106: was 3442.25
$ws.Cells.Item(106, 9).Value = 0  # I106: was 3444
$ws.Cells.Item(106, 10).Value = 3406.25  # J106: was 3441.6667
$ws.Cells.Item(106, 11).Value = 0  # K106: was 10332
$ws.Cells.Item(106, 12).Value = 10218.75  # L106: was 10325.0001
$ws.Cells.Item(106, 13).ClearContents()  # M106: was -9386
$ws.Cells.Item(106, 14).Value = -12110.75  # N106: was -12217.0001
$ws.Cells.Item(107, 8).Value = 1387  # H107: was 1427
$ws.Cells.Item(107, 9).Value = 1298.25  # I107: was 1238.6
$ws.Cells.Item(107, 10).Value = 1446.1666  # J107: was 1615.4
$ws.Cells.Item(107, 11).Value = 3894.75  # K107: was 3715.8
$ws.Cells.Item(107, 12).Value = 4338.4998  # L107: was 4846.200000000001
$ws.Cells.Item(107, 13).Value = -1974.75  # M107: was -1795.8
$ws.Cells.Item(107, 14).Value = -8178.4998  # N107: was -8686.200000000001
$ws.Cells.Item(120, 8).Value = 15919.167  # H120: was 17403
$ws.Cells.Item(120, 9).Value = 11200  # I120: was 13000
$ws.Cells.Item(120, 11).Value = 33600  # K120: was 39000
$ws.Cells.Item(120, 13).Value = -28762  # M120: was -34162
$ws.Cells.Item(132, 8).Value = 1705  # H132: was 1906.25
$ws.Cells.Item(132, 10).Value = 1966.6666  # J132: was 2500
$ws.Cells.Item(132, 12).Value = 17699.9994  # L132: was 22500
$ws.Cells.Item(132, 14).Value = -22759.9994  # N132: was -27560
$ws.Cells.Item(140, 8).Value = 1527.85  # H140: was 1605.1052
$ws.Cells.Item(140, 9).Value = 1198.6471  # I140: was 1269.8125
$ws.Cells.Item(140, 11).Value = 3595.9413  # K140: was 3809.4375
$ws.Cells.Item(140, 13).Value = 1584.0587  # M140: was 1370.5625

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 8489.666999999999  # H2: was 9256.909
$ws.Cells.Item(2, 9).Value = 80  # I2: was 83.333336
$ws.Cells.Item(2, 10).Value = 12694.5  # J2: was 12697
$ws.Cells.Item(2, 11).Value = 80  # K2: was 83.333336
$ws.Cells.Item(2, 12).Value = 12694.5  # L2: was 12697
$ws.Cells.Item(2, 13).Value = 33  # M2: was 29.666664
$ws.Cells.Item(2, 14).Value = -12920.5  # N2: was -12923
$ws.Cells.Item(17, 8).Value = 2401.0344  # H17: was 3213.5483
$ws.Cells.Item(17, 10).Value = 33952.5  # J17: was 24473.75
$ws.Cells.Item(17, 12).Value = 33952.5  # L17: was 24473.75
$ws.Cells.Item(17, 14).Value = -34288.5  # N17: was -24809.75
$ws.Cells.Item(32, 8).Value = 0  # H32: was 17000000
$ws.Cells.Item(32, 10).Value = 0  # J32: was 17000000
$ws.Cells.Item(32, 12).Value = 0  # L32: was 17000000
$ws.Cells.Item(32, 14).ClearContents()  # N32: was -17000592
$ws.Cells.Item(41, 8).Value = 12446.143  # H41: was 11015.375
$ws.Cells.Item(41, 9).Value = 12446.143  # I41: was 11015.375
$ws.Cells.Item(41, 11).Value = 12446.143  # K41: was 11015.375
$ws.Cells.Item(41, 13).Value = -12091.143  # M41: was -10660.375
$ws.Cells.Item(122, 8).Value = 1818.2858  # H122: was 1772.75
$ws.Cells.Item(122, 9).Value = 1781  # I122: was 1726.5
$ws.Cells.Item(122, 11).Value = 5343  # K122: was 5179.5
$ws.Cells.Item(122, 13).Value = -2893  # M122: was -2729.5
$ws.Cells.Item(126, 8).Value = 2812.318  # H126: was 2893.55
$ws.Cells.Item(126, 9).Value = 1760.1666  # I126: was 1783.8182
$ws.Cells.Item(126, 10).Value = 4074.9  # J126: was 4249.8887
$ws.Cells.Item(126, 11).Value = 5280.4998  # K126: was 5351.4546
$ws.Cells.Item(126, 12).Value = 12224.7  # L126: was 12749.6661
$ws.Cells.Item(126, 13).Value = -2810.4998  # M126: was -2881.4546
$ws.Cells.Item(126, 14).Value = -17164.7  # N126: was -17689.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6760  # H7: was 7942
$ws.Cells.Item(7, 9).Value = 815  # I7: was 803.3333
$ws.Cells.Item(7, 11).Value = 815  # K7: was 803.3333
$ws.Cells.Item(7, 13).Value = -703  # M7: was -691.3333
$ws.Cells.Item(23, 8).Value = 9989  # H23: was 7494
$ws.Cells.Item(23, 9).Value = 0  # I23: was 4999
$ws.Cells.Item(23, 11).Value = 0  # K23: was 4999
$ws.Cells.Item(23, 13).ClearContents()  # M23: was -4769
$ws.Cells.Item(101, 8).Value = 23500  # H101: was 64489.5
$ws.Cells.Item(101, 10).Value = 23500  # J101: was 64489.5
$ws.Cells.Item(101, 12).Value = 23500  # L101: was 64489.5
$ws.Cells.Item(101, 14).Value = -29990  # N101: was -70979.5
$ws.Cells.Item(111, 8).Value = 0  # H111: was 98385
$ws.Cells.Item(111, 10).Value = 0  # J111: was 98385
$ws.Cells.Item(111, 12).Value = 0  # L111: was 98385
$ws.Cells.Item(111, 14).ClearContents()  # N111: was -106565
$ws.Cells.Item(126, 8).Value = 6760  # H126: was 7942
$ws.Cells.Item(126, 9).Value = 815  # I126: was 803.3333
$ws.Cells.Item(126, 11).Value = 2445  # K126: was 2409.9999
$ws.Cells.Item(126, 13).Value = 25  # M126: was 60.0001000000002
$ws.Cells.Item(128, 8).Value = 55998  # H128: was 63140.285
$ws.Cells.Item(128, 10).Value = 55998  # J128: was 63140.285
$ws.Cells.Item(128, 12).Value = 55998  # L128: was 63140.285
$ws.Cells.Item(128, 14).Value = -65958  # N128: was -73100.285

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(106, 8).Value = 40985.4  # H106: was 41652
$ws.Cells.Item(106, 10).Value = 40985.4  # J106: was 41652
$ws.Cells.Item(106, 12).Value = 40985.4  # L106: was 41652
$ws.Cells.Item(106, 14).Value = -43509.4  # N106: was -44176
$ws.Cells.Item(108, 8).Value = 0  # H108: was 50000
$ws.Cells.Item(108, 9).Value = 0  # I108: was 50000
$ws.Cells.Item(108, 11).Value = 0  # K108: was 50000
$ws.Cells.Item(108, 13).ClearContents()  # M108: was -46160
